$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.798.46'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '2.313.33'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'112.46"
$ws.Range("E5").Value = '  +17.49%  '
$ws.Range("D6").Value = "'271.80"
$ws.Range("E6").Value = '  +1.18%  '
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = "'0.627"
$ws.Range("E9").Value = '  +1.81%  '
$ws.Range("D10").Value = "'47.45"
$ws.Range("E10").Value = '  +6.77%  '
$ws.Range("D11").Value = "'0.0943"
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("D12").Value = "'8.81"
$ws.Range("E12").Value = '  +11.38%  '
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("D14").Value = "'15.81"
$ws.Range("E14").Value = '  +3.18%  '
$ws.Range("D15").Value = '2.658.42'
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").Value = "'0.859"
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("D17").Value = '2.324.44'
$ws.Range("E17").Value = '  +0.78%  '
$ws.Range("D18").Value = '43.865.78'
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").Value = "'0.0000110"
$ws.Range("E19").Value = '  +1.86%  '
$ws.Range("D20").Value = "'6.57"
$ws.Range("E20").Value = '  +4.44%  '
$ws.Range("D21").Value = "'72.67"
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").Value = "'2.53"
$ws.Range("E22").Value = '  +7.81%  '
$ws.Range("D23").Value = "'234.81"
$ws.Range("E23").Value = '  -1.27%  '
$ws.Range("D24").Value = "'9.61"
$ws.Range("E24").Value = '  +5.52%  '
$ws.Range("D25").Value = "'2.92"
$ws.Range("E25").Value = '  +16.64%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").Value = "'11.52"
$ws.Range("E27").Value = '  +1.77%  '
$ws.Range("D28").Value = "'43.23"
$ws.Range("E28").Value = '  +13.83%  '
$ws.Range("E29").Value = '  -0.53%  '
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("D31").Value = "'178.39"
$ws.Range("E31").Value = '  +1.99%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = "'0.0947"
$ws.Range("E32").Value = '  +6.15%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = "'21.86"
$ws.Range("E33").Value = '  -1.49%  '
$ws.Range("D34").Value = "'5.63"
$ws.Range("E34").Value = '  +3.29%  '
$ws.Range("D35").Value = "'4.83"
$ws.Range("E35").Value = '  +8.15%  '
$ws.Range("D36").Value = "'0.128"
$ws.Range("E36").Value = '  +1.08%  '
$ws.Range("D37").Value = "'0.113"
$ws.Range("E37").Value = '  +3.95%  '
$ws.Range("D38").Value = "'3.96"
$ws.Range("E38").Value = '  +20.89%  '
$ws.Range("D39").Value = "'0.0359"
$ws.Range("E39").Value = '  +0.24%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = "'0.243"
$ws.Range("E40").Value = '  +1.75%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").Value = "'2.40"
$ws.Range("E41").Value = '  +1.15%  '
$ws.Range("D42").Value = "'70.29"
$ws.Range("E42").Value = '  +12.56%  '
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("D44").Value = "'1.39"
$ws.Range("E44").Value = '  +1.21%  '
$ws.Range("D45").Value = "'12.38"
$ws.Range("E45").Value = '  +2.31%  '
$ws.Range("D46").Value = "'5.57"
$ws.Range("E46").Value = '  +5.36%  '
$ws.Range("D47").Value = "'8.83"
$ws.Range("E47").Value = '  -2.39%  '
$ws.Range("D48").Value = "'0.101"
$ws.Range("E48").Value = '  -1.51%  '
$ws.Range("D49").Value = "'100.07"
$ws.Range("E49").Value = '  -0.23%  '
$ws.Range("E50").Value = '  +2.10%  '
$ws.Range("D51").Value = "'0.461"
$ws.Range("E51").Value = '  +9.69%  '
